$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header shared string: "value" -> "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# Update / extend data rows 2-84 (dates in col A, values in col B)
$ws.Range("A2").Value = 38398
$ws.Range("B2").Value = -0.7
$ws.Range("A3").Value = 38487
$ws.Range("B3").Value = 0.4
$ws.Range("A4").Value = 38579
$ws.Range("B4").Value = -0.2
$ws.Range("A5").Value = 38671
$ws.Range("B5").Value = 0.7
$ws.Range("A6").Value = 38763
$ws.Range("B6").Value = 0
$ws.Range("A7").Value = 38852
$ws.Range("B7").Value = -0.4
$ws.Range("A8").Value = 38944
$ws.Range("B8").Value = 0.1
$ws.Range("A9").Value = 39036
$ws.Range("B9").Value = -1.6
$ws.Range("A10").Value = 39128
$ws.Range("B10").Value = 2
$ws.Range("A11").Value = 39217
$ws.Range("B11").Value = -1.1
$ws.Range("A12").Value = 39309
$ws.Range("B12").Value = 0.4
$ws.Range("A13").Value = 39401
$ws.Range("B13").Value = 0
$ws.Range("A14").Value = 39493
$ws.Range("B14").Value = 0.7
$ws.Range("A15").Value = 39583
$ws.Range("B15").Value = -0.4
$ws.Range("A16").Value = 39675
$ws.Range("B16").Value = 0.9
$ws.Range("A17").Value = 39767
$ws.Range("B17").Value = 0.6
$ws.Range("A18").Value = 39859
$ws.Range("B18").Value = -0.5
$ws.Range("A19").Value = 39948
$ws.Range("B19").Value = -2
$ws.Range("A20").Value = 40040
$ws.Range("B20").Value = 1.5
$ws.Range("A21").Value = 40132
$ws.Range("B21").Value = -1.1
$ws.Range("A22").Value = 40224
$ws.Range("B22").Value = 1.9
$ws.Range("A23").Value = 40313
$ws.Range("B23").Value = 0.1
$ws.Range("A24").Value = 40405
$ws.Range("B24").Value = -0.3
$ws.Range("A25").Value = 40497
$ws.Range("B25").Value = -0.5
$ws.Range("A26").Value = 40589
$ws.Range("B26").Value = -0.4
$ws.Range("A27").Value = 40678
$ws.Range("B27").Value = 0.3
$ws.Range("A28").Value = 40770
$ws.Range("B28").Value = -0.4
$ws.Range("A29").Value = 40862
$ws.Range("B29").Value = 0
$ws.Range("A30").Value = 40954
$ws.Range("B30").Value = -0.4
$ws.Range("A31").Value = 41044
$ws.Range("B31").Value = 0
$ws.Range("A32").Value = 41136
$ws.Range("B32").Value = -0.3
$ws.Range("A33").Value = 41228
$ws.Range("B33").Value = 0.4
$ws.Range("A34").Value = 41320
$ws.Range("B34").Value = -0.1
$ws.Range("A35").Value = 41409
$ws.Range("B35").Value = -0.3
$ws.Range("A36").Value = 41501
$ws.Range("B36").Value = 0.2
$ws.Range("A37").Value = 41593
$ws.Range("B37").Value = -0.2
$ws.Range("A38").Value = 41685
$ws.Range("B38").Value = 0.7000000000000001
$ws.Range("A39").Value = 41774
$ws.Range("B39").Value = 0.1
$ws.Range("A40").Value = 41866
$ws.Range("B40").Value = -0.5
$ws.Range("A41").Value = 41958
$ws.Range("B41").Value = 0.4
$ws.Range("A42").Value = 42050
$ws.Range("B42").Value = -0.3
$ws.Range("A43").Value = 42139
$ws.Range("B43").Value = -0.3
$ws.Range("A44").Value = 42231
$ws.Range("B44").Value = 0.2
$ws.Range("A45").Value = 42323
$ws.Range("B45").Value = 0.1
$ws.Range("A46").Value = 42415
$ws.Range("B46").Value = 0.1
$ws.Range("A47").Value = 42505
$ws.Range("B47").Value = -0.2
$ws.Range("A48").Value = 42597
$ws.Range("B48").Value = 0
$ws.Range("A49").Value = 42689
$ws.Range("B49").Value = 0.4
$ws.Range("A50").Value = 42781
$ws.Range("B50").Value = -0.4
$ws.Range("A51").Value = 42870
$ws.Range("B51").Value = 0.2
$ws.Range("A52").Value = 42962
$ws.Range("B52").Value = 0.4
$ws.Range("A53").Value = 43054
$ws.Range("B53").Value = 0
$ws.Range("A54").Value = 43146
$ws.Range("B54").Value = -0.1
$ws.Range("A55").Value = 43235
$ws.Range("B55").Value = 0.4
$ws.Range("A56").Value = 43327
$ws.Range("B56").Value = 0.7
$ws.Range("A57").Value = 43419
$ws.Range("B57").Value = -0.6
$ws.Range("A58").Value = 43511
$ws.Range("B58").Value = -0.6
$ws.Range("A59").Value = 43600
$ws.Range("B59").Value = 0.3
$ws.Range("A60").Value = 43692
$ws.Range("B60").Value = -0.7
$ws.Range("A61").Value = 43784
$ws.Range("B61").Value = 0.6
$ws.Range("A62").Value = 43876
$ws.Range("B62").Value = 0.3
$ws.Range("A63").Value = 43966
$ws.Range("B63").Value = 0.3
$ws.Range("A64").Value = 44058
$ws.Range("B64").Value = -0.4440571223929872
$ws.Range("A65").Value = 44150
$ws.Range("B65").Value = -0.7255945204468831
$ws.Range("A66").Value = 44242
$ws.Range("B66").Value = -0.5292660609007143
$ws.Range("A67").Value = 44331
$ws.Range("B67").Value = -0.1550786956675604
$ws.Range("A68").Value = 44423
$ws.Range("B68").Value = -2.168330733759602
$ws.Range("A69").Value = 44515
$ws.Range("B69").Value = -0.03982694963614287
$ws.Range("A70").Value = 44607
$ws.Range("B70").Value = 0.2669401745841223
$ws.Range("A71").Value = 44696
$ws.Range("B71").Value = 0.03791487406588956
$ws.Range("A72").Value = 44788
$ws.Range("B72").Value = -0.04567208272808071
$ws.Range("A73").Value = 44880
$ws.Range("B73").Value = -0.5154625125417773
$ws.Range("A74").Value = 44972
$ws.Range("B74").Value = -0.1813602613933202
$ws.Range("A75").Value = 45061
$ws.Range("B75").Value = -0.01480819732384536
$ws.Range("A76").Value = 45153
$ws.Range("B76").Value = 0.02918400950819283
$ws.Range("A77").Value = 45245
$ws.Range("B77").Value = -0.03321544329283629
$ws.Range("A78").Value = 45337
$ws.Range("B78").Value = 0.00001303303454188581
$ws.Range("A79").Value = 45427
$ws.Range("B79").Value = -0.006125572440376981
$ws.Range("A80").Value = 45519
$ws.Range("B80").Value = 0.04879937325030748
$ws.Range("A81").Value = 45611
$ws.Range("B81").Value = 0.0477695913607396
$ws.Range("A82").Value = 45703
$ws.Range("B82").Value = 0.4714513528429705
$ws.Range("A83").Value = 45792
$ws.Range("B83").Value = -0.02605454389395597
$ws.Range("A84").Value = 45884
$ws.Range("B84").Value = 0.04549112474043772

# New rows (54-84) need the same date-cell style/number-format as the existing rows
$ws.Range("A2").Copy()
$ws.Range("A54:A84").PasteSpecial(-4122)
$excel.CutCopyMode = 0
